$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. D63: fill in the missing "idPagamento" (was blank, now has a value) ---
# Column D stores payment IDs as text (e.g. other rows like D2 = "76158393297"),
# so force a Text number format before assigning a purely-numeric string —
# otherwise Excel would auto-coerce it to a Number.
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value2 = "77481297146"

# --- 2. Append a new record as row 69 for "Vitor Ito" ---
$ws.Range("A69").Value2 = "Vitor Ito"
$ws.Range("B69").Value2 = 1578424633
$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value2 = "11988776655"
# D69 (idPagamento) is left blank, matching the source row's empty value.
$ws.Range("E69").Value2 = 1
$ws.Range("F69").Value2 = 2
$ws.Range("G69").Value2 = 3
$ws.Range("H69").Value2 = 4
$ws.Range("I69").Value2 = 5
$ws.Range("J69").Value2 = 9
$ws.Range("K69").Value2 = 10
$ws.Range("L69").Value2 = 22
$ws.Range("M69").Value2 = 28
$ws.Range("N69").Value2 = 29
$ws.Range("O69").Value2 = "Não"

# --- 3. Make sure the "number stored as text" warning is suppressed over the
#        full data range, same as it was for the previous A1:O68 extent. ---
$full = $ws.Range("A1:O69")
try {
    $full.Errors.Item(7).Ignore = $true
} catch {
}
